# Update timelog with today's meeting entry (row 8 of the timesheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the next timelog row: date, hours worked, and task description.
$ws.Range("A8").Value = "16/9/14"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "Combined Requirements Documents for the computer and pushed them into repository"
